$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Thbs1"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 122.253015
$ws.Range("H2").Value = 366.759045
$ws.Range("I2").Value = 0.1988639364328829
$ws.Range("J2").Value = 0.1988639364328829
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.481489333333333
$ws.Range("N2").Value = 7.444467999999999
$ws.Range("O2").Value = 0.2345069082418988
$ws.Range("P2").Value = 0.2345069082418987
$ws.Range("Q2").Value = 303.36955269034
$ws.Range("R2").Value = 2730.325974213059
$ws.Range("S2").Value = 0.04663496689368886
$ws.Range("T2").Value = 0.04663496689368884

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Thbs1"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 122.253015
$ws.Range("H3").Value = 366.759045
$ws.Range("I3").Value = 0.1988639364328829
$ws.Range("J3").Value = 0.1988639364328829
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.245227
$ws.Range("N3").Value = 21.735681
$ws.Range("O3").Value = 0.6846919551326144
$ws.Range("P3").Value = 0.6846919551326142
$ws.Range("Q3").Value = 885.750845109405
$ws.Range("R3").Value = 7971.757605984645
$ws.Range("S3").Value = 0.1361605374415985
$ws.Range("T3").Value = 0.1361605374415985

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Thbs1"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 122.253015
$ws.Range("H4").Value = 366.759045
$ws.Range("I4").Value = 0.1988639364328829
$ws.Range("J4").Value = 0.1988639364328829
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2001876666666667
$ws.Range("N4").Value = 0.600563
$ws.Range("O4").Value = 0.01891823194544989
$ws.Range("P4").Value = 0.01891823194544989
$ws.Range("Q4").Value = 24.473545815815
$ws.Range("R4").Value = 220.261912342335
$ws.Range("S4").Value = 0.003762154075022482
$ws.Range("T4").Value = 0.003762154075022481

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Thbs1"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 122.253015
$ws.Range("H5").Value = 366.759045
$ws.Range("I5").Value = 0.1988639364328829
$ws.Range("J5").Value = 0.1988639364328829
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6548283333333333
$ws.Range("N5").Value = 1.964485
$ws.Range("O5").Value = 0.06188290468003712
$ws.Range("P5").Value = 0.06188290468003711
$ws.Range("Q5").Value = 80.054738057425
$ws.Range("R5").Value = 720.4926425168251
$ws.Range("S5").Value = 0.01230627802257305
$ws.Range("T5").Value = 0.01230627802257305

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Thbs1"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 132.5447616666667
$ws.Range("H6").Value = 397.634285
$ws.Range("I6").Value = 0.2156050961899926
$ws.Range("J6").Value = 0.2156050961899926
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.481489333333333
$ws.Range("N6").Value = 7.444467999999999
$ws.Range("O6").Value = 0.2345069082418988
$ws.Range("P6").Value = 0.2345069082418987
$ws.Range("Q6").Value = 328.9084122650421
$ws.Range("R6").Value = 2960.175710385379
$ws.Range("S6").Value = 0.05056088450871235
$ws.Range("T6").Value = 0.05056088450871234

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Thbs1"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 132.5447616666667
$ws.Range("H7").Value = 397.634285
$ws.Range("I7").Value = 0.2156050961899926
$ws.Range("J7").Value = 0.2156050961899926
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.245227
$ws.Range("N7").Value = 21.735681
$ws.Range("O7").Value = 0.6846919551326144
$ws.Range("P7").Value = 0.6846919551326142
$ws.Range("Q7").Value = 960.3168859358983
$ws.Range("R7").Value = 8642.851973423085
$ws.Range("S7").Value = 0.1476230748468814
$ws.Range("T7").Value = 0.1476230748468814

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Thbs1"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 132.5447616666667
$ws.Range("H8").Value = 397.634285
$ws.Range("I8").Value = 0.2156050961899926
$ws.Range("J8").Value = 0.2156050961899926
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2001876666666667
$ws.Range("N8").Value = 0.600563
$ws.Range("O8").Value = 0.01891823194544989
$ws.Range("P8").Value = 0.01891823194544989
$ws.Range("Q8").Value = 26.53382656693944
$ws.Range("R8").Value = 238.804439102455
$ws.Range("S8").Value = 0.004078867218343315
$ws.Range("T8").Value = 0.004078867218343314

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Thbs1"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 132.5447616666667
$ws.Range("H9").Value = 397.634285
$ws.Range("I9").Value = 0.2156050961899926
$ws.Range("J9").Value = 0.2156050961899926
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6548283333333333
$ws.Range("N9").Value = 1.964485
$ws.Range("O9").Value = 0.06188290468003712
$ws.Range("P9").Value = 0.06188290468003711
$ws.Range("Q9").Value = 86.79406537424723
$ws.Range("R9").Value = 781.146588368225
$ws.Range("S9").Value = 0.01334226961605555
$ws.Range("T9").Value = 0.01334226961605554

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Thbs1"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 320.0894206666666
$ws.Range("H10").Value = 960.2682619999999
$ws.Range("I10").Value = 0.5206762565675317
$ws.Range("J10").Value = 0.5206762565675317
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.481489333333333
$ws.Range("N10").Value = 7.444467999999999
$ws.Range("O10").Value = 0.2345069082418988
$ws.Range("P10").Value = 0.2345069082418987
$ws.Range("Q10").Value = 794.2984830971793
$ws.Range("R10").Value = 7148.686347874615
$ws.Range("S10").Value = 0.1221021791226175
$ws.Range("T10").Value = 0.1221021791226175

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Thbs1"
$ws.Range("C11").Value = "Itgb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 320.0894206666666
$ws.Range("H11").Value = 960.2682619999999
$ws.Range("I11").Value = 0.5206762565675317
$ws.Range("J11").Value = 0.5206762565675317
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 7.245227
$ws.Range("N11").Value = 21.735681
$ws.Range("O11").Value = 0.6846919551326144
$ws.Range("P11").Value = 0.6846919551326142
$ws.Range("Q11").Value = 2319.120513028491
$ws.Range("R11").Value = 20872.08461725642
$ws.Range("S11").Value = 0.356502844100354
$ws.Range("T11").Value = 0.3565028441003539

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Thbs1"
$ws.Range("C12").Value = "Itgb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 320.0894206666666
$ws.Range("H12").Value = 960.2682619999999
$ws.Range("I12").Value = 0.5206762565675317
$ws.Range("J12").Value = 0.5206762565675317
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2001876666666667
$ws.Range("N12").Value = 0.600563
$ws.Range("O12").Value = 0.01891823194544989
$ws.Range("P12").Value = 0.01891823194544989
$ws.Range("Q12").Value = 64.0779542479451
$ws.Range("R12").Value = 576.7015882315059
$ws.Range("S12").Value = 0.009850274190233143
$ws.Range("T12").Value = 0.009850274190233141

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Thbs1"
$ws.Range("C13").Value = "Itgb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 320.0894206666666
$ws.Range("H13").Value = 960.2682619999999
$ws.Range("I13").Value = 0.5206762565675317
$ws.Range("J13").Value = 0.5206762565675317
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6548283333333333
$ws.Range("N13").Value = 1.964485
$ws.Range("O13").Value = 0.06188290468003712
$ws.Range("P13").Value = 0.06188290468003711
$ws.Range("Q13").Value = 209.6036218527855
$ws.Range("R13").Value = 1886.43259667507
$ws.Range("S13").Value = 0.03222095915432711
$ws.Range("T13").Value = 0.03222095915432711

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Thbs1"
$ws.Range("C14").Value = "Itgb3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 39.86989333333333
$ws.Range("H14").Value = 119.60968
$ws.Range("I14").Value = 0.06485471080959287
$ws.Range("J14").Value = 0.06485471080959287
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.481489333333333
$ws.Range("N14").Value = 7.444467999999999
$ws.Range("O14").Value = 0.2345069082418988
$ws.Range("P14").Value = 0.2345069082418987
$ws.Range("Q14").Value = 98.93671502780441
$ws.Range("R14").Value = 890.4304352502398
$ws.Range("S14").Value = 0.01520887771688008
$ws.Range("T14").Value = 0.01520887771688007

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Thbs1"
$ws.Range("C15").Value = "Itgb3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 39.86989333333333
$ws.Range("H15").Value = 119.60968
$ws.Range("I15").Value = 0.06485471080959287
$ws.Range("J15").Value = 0.06485471080959287
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 7.245227
$ws.Range("N15").Value = 21.735681
$ws.Range("O15").Value = 0.6846919551326144
$ws.Range("P15").Value = 0.6846919551326142
$ws.Range("Q15").Value = 288.8664276657867
$ws.Range("R15").Value = 2599.79784899208
$ws.Range("S15").Value = 0.04440549874378045
$ws.Range("T15").Value = 0.04440549874378043

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Thbs1"
$ws.Range("C16").Value = "Itgb3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 39.86989333333333
$ws.Range("H16").Value = 119.60968
$ws.Range("I16").Value = 0.06485471080959287
$ws.Range("J16").Value = 0.06485471080959287
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.2001876666666667
$ws.Range("N16").Value = 0.600563
$ws.Range("O16").Value = 0.01891823194544989
$ws.Range("P16").Value = 0.01891823194544989
$ws.Range("Q16").Value = 7.981460916648888
$ws.Range("R16").Value = 71.83314824984
$ws.Range("S16").Value = 0.001226936461850954
$ws.Range("T16").Value = 0.001226936461850954

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Thbs1"
$ws.Range("C17").Value = "Itgb3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 39.86989333333333
$ws.Range("H17").Value = 119.60968
$ws.Range("I17").Value = 0.06485471080959287
$ws.Range("J17").Value = 0.06485471080959287
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6548283333333333
$ws.Range("N17").Value = 1.964485
$ws.Range("O17").Value = 0.06188290468003712
$ws.Range("P17").Value = 0.06188290468003711
$ws.Range("Q17").Value = 26.10793580164444
$ws.Range("R17").Value = 234.9714222148
$ws.Range("S17").Value = 0.004013397887081408
$ws.Range("T17").Value = 0.004013397887081408
